# Generate Report for Handoff
# Row 3 (the "b.md" source file) moves from "Handed back" to "Ready for
# handoff" status, with a freshly generated handoff file + handoff
# datetime for each target language sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet -------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-10 00:35:01"
foreach ($hl in $wsZh.Hyperlinks) {
    $r = $hl.Range
    if ($r.Row -eq 3 -and $r.Column -eq 3) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-10 00:35:09"
foreach ($hl in $wsDe.Hyperlinks) {
    $r = $hl.Range
    if ($r.Row -eq 3 -and $r.Column -eq 3) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
